$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename "CSS in JavaScript" (row 8, column C) to "CSS in JavaScript (CSS-in-JS)"
#    and give it the "typed in" default-font style (matches the new font/cellXf
#    that appears in styles.xml for this cell).
$c8 = $ws.Range("C8")
$c8.Value = "CSS in JavaScript (CSS-in-JS)"
$c8.Font.ThemeFont = 1

# 2. Update the title-generating formulas so "Learning" moves from before the
#    topic name to after it:
#    ="<title>Learning "&C2&" Resources - ..." -> ="<title>"&C2&" Learning Resources - ..."
$ws.Range("E2").Formula = '="<title>"&C2&" Learning Resources - Front-End Developer Learning Roadmap</title>"'
$ws.Range("E3:E59").Formula = '="<title>"&C3&" Learning Resources - Front-End Developer Learning Roadmap</title>"'

# 3. Update the view/selection state to match the saved workbook
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("E13").Select() | Out-Null
